# Actualización automática 2025-12-07 17:44:30
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 12 (BUENAÑO VITERI MARJORIE LETICIA)
$wsGrupo.Range("D12").Value = 915.84
$wsGrupo.Range("M12").Value = 12993.77

# Row 60 totals (counts of "X de 58")
$wsGrupo.Range("D60").Value = "2 de 58"
$wsGrupo.Range("M60").Value = "3 de 58"

# --- Sheet 2: "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("G7").Value = 1000
$wsMensual.Range("G10").Value = 10000

$wsMensual.Range("F12").Value = 13909.61
$wsMensual.Range("G12").Value = 14000

$wsMensual.Range("G15").Value = 0
$wsMensual.Range("G16").Value = 1000
$wsMensual.Range("G19").Value = 0
$wsMensual.Range("G22").Value = 500
$wsMensual.Range("G26").Value = 2500
$wsMensual.Range("G27").Value = 2000
$wsMensual.Range("G29").Value = 1000
$wsMensual.Range("G30").Value = 2000
$wsMensual.Range("G41").Value = 500
$wsMensual.Range("G46").Value = 1500
$wsMensual.Range("G47").Value = 2000
$wsMensual.Range("G53").Value = 2500
$wsMensual.Range("G58").Value = 1000
$wsMensual.Range("G59").Value = 0

# Totals row
$wsMensual.Range("F60").Value = 18460.28
$wsMensual.Range("G60").Value = 48000
